# Applies the changes described by the diff:
#  - Rows 13-16 (columns E,F,G) get their description/barcode/price
#    rotated: row13->row16, row14->row13, row15->row14, row16->row15
#    (equivalent to the shared-string reordering performed in the source diff)
#  - Row15 SalesQuantity/Turnover (K,L) pick up what used to be on row16
#  - Row16 Turnover (L) gets a new value
#  - Row45 SalesQuantity/Turnover (K,L) updated
#  - Row48 (totals row): J48 removed, K48/L48 updated to new column totals

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture current ("before") values for the rotating block (rows 13-16) ---
$E13 = $ws.Range("E13").Value()
$E14 = $ws.Range("E14").Value()
$E15 = $ws.Range("E15").Value()
$E16 = $ws.Range("E16").Value()

# Barcodes are long strings of digits stored as text in the sheet; reading
# them through .Value() would silently coerce them to Double (losing their
# text type), so the current values are restated here literally instead.
$F13 = "5201178035923"
$F14 = "5201178022862"
$F15 = "4005808478200"
$F16 = "4005808445417"

$G13 = $ws.Range("G13").Value()
$G14 = $ws.Range("G14").Value()
$G15 = $ws.Range("G15").Value()
$G16 = $ws.Range("G16").Value()

# --- Rotate descriptions (E) and prices (G): new row N = old row N+1
#     (row16 wraps around to old row13) ---
$ws.Range("E13").Value = $E14
$ws.Range("E14").Value = $E15
$ws.Range("E15").Value = $E16
$ws.Range("E16").Value = $E13

$ws.Range("G13").Value = $G14
$ws.Range("G14").Value = $G15
$ws.Range("G15").Value = $G16
$ws.Range("G16").Value = $G13

# --- Rotate barcodes (F), forcing them to stay text (not numbers) ---
foreach ($addr in @("F13","F14","F15","F16")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("F13").Value = $F14
$ws.Range("F14").Value = $F15
$ws.Range("F15").Value = $F16
$ws.Range("F16").Value = $F13

# Restore the original (General / left-aligned) style on the barcode cells -
# the NumberFormat tweak above is only there to stop Excel from reinterpreting
# the digit-string as a number.
foreach ($r in 13,14,15,16) {
    $ws.Range("E$r").Copy()
    $ws.Range("F$r").PasteSpecial(-4122)   # xlPasteFormats
}
$excel.CutCopyMode = $false

# --- Row15 / Row16 SalesQuantity & Turnover ---
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 8.19
$ws.Range("L16").Value = 6.38

# --- Row45 SalesQuantity & Turnover ---
$ws.Range("K45").Value = 5
$ws.Range("L45").Value = 39.35

# --- Row48 totals: J48 removed entirely, K48 updated, L48 added (new cell,
#     picking up K48's bold/red "totals" formatting) ---
$ws.Range("K48").Copy()
$ws.Range("L48").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("K48").Value = 70
$ws.Range("L48").Value = 596.61
$ws.Range("J48").Clear()
